# Apply cryptocurrency price/volume updates to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='28.074.56'},
    @{Cell='E2'; Value='  -0.50%  '},
    @{Cell='D3'; Value='1.799.88'},
    @{Cell='E3'; Value='  -0.22%  '},
    @{Cell='D4'; Value='1.003'},
    @{Cell='E4'; Value='  -0.24%  '},
    @{Cell='D5'; Value='311.23'},
    @{Cell='E5'; Value='  -1.22%  '},
    @{Cell='D7'; Value='0.5096'},
    @{Cell='E7'; Value='  -2.55%  '},
    @{Cell='D8'; Value='0.3861'},
    @{Cell='E8'; Value='  +1.19%  '},
    @{Cell='D9'; Value='0.07722'},
    @{Cell='E9'; Value='  -2.58%  '},
    @{Cell='E10'; Value='  +0.14%  '},
    @{Cell='E11'; Value='  -2.03%  '},
    @{Cell='E12'; Value='  -0.37%  '},
    @{Cell='E13'; Value='  -0.33%  '},
    @{Cell='D14'; Value='20.33'},
    @{Cell='E14'; Value='  -1.62%  '},
    @{Cell='B15'; Value='Chainlink'},
    @{Cell='C15'; Value='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'},
    @{Cell='D15'; Value='7.279'},
    @{Cell='E15'; Value='  -0.96%  '},
    @{Cell='B16'; Value='WrappedEther'},
    @{Cell='C16'; Value='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'},
    @{Cell='D16'; Value='1.793.91'},
    @{Cell='E16'; Value='  -0.56%  '},
    @{Cell='D17'; Value='92.10'},
    @{Cell='E17'; Value='  -0.57%  '},
    @{Cell='D18'; Value='0.00001074'},
    @{Cell='E18'; Value='  -1.34%  '},
    @{Cell='D19'; Value='0.06554'},
    @{Cell='E19'; Value='  -0.60%  '},
    @{Cell='D20'; Value='1.002'},
    @{Cell='E20'; Value='  -0.26%  '},
    @{Cell='D21'; Value='17.23'},
    @{Cell='E21'; Value='  -1.26%  '},
    @{Cell='D22'; Value='5.949'},
    @{Cell='E22'; Value='  -0.28%  '},
    @{Cell='D23'; Value='28.107.58'},
    @{Cell='E23'; Value='  -0.61%  '},
    @{Cell='E24'; Value='  -0.54%  '},
    @{Cell='D25'; Value='2.247'},
    @{Cell='E25'; Value='  +0.63%  '},
    @{Cell='D26'; Value='160.68'},
    @{Cell='E26'; Value='  +1.99%  '},
    @{Cell='D27'; Value='2.414'},
    @{Cell='E27'; Value='  +0.45%  '},
    @{Cell='D28'; Value='2.005.65'},
    @{Cell='E28'; Value='  -0.22%  '},
    @{Cell='E29'; Value='  -1.36%  '},
    @{Cell='D30'; Value='127.09'},
    @{Cell='E30'; Value='  +3.20%  '},
    @{Cell='D31'; Value='0.1086'},
    @{Cell='E31'; Value='  -1.85%  '},
    @{Cell='E32'; Value='  -1.61%  '},
    @{Cell='D33'; Value='3.649'},
    @{Cell='E33'; Value='  -0.53%  '},
    @{Cell='D34'; Value='5.531'},
    @{Cell='E34'; Value='  -0.69%  '},
    @{Cell='D35'; Value='0.07028'},
    @{Cell='E35'; Value='  -2.37%  '},
    @{Cell='D36'; Value='9.035'},
    @{Cell='E36'; Value='  +2.67%  '},
    @{Cell='D37'; Value='0.02345'},
    @{Cell='E37'; Value='  +1.23%  '},
    @{Cell='D38'; Value='0.2164'},
    @{Cell='E38'; Value='  -0.36%  '},
    @{Cell='D39'; Value='5.024'},
    @{Cell='E39'; Value='  -0.23%  '},
    @{Cell='E40'; Value='  -6.04%  '},
    @{Cell='D41'; Value='0.6115'},
    @{Cell='E41'; Value='  -1.55%  '},
    @{Cell='D43'; Value='1.151'},
    @{Cell='E43'; Value='  -1.24%  '},
    @{Cell='D44'; Value='13.19'},
    @{Cell='E44'; Value='  -0.30%  '},
    @{Cell='D45'; Value='0.5912'},
    @{Cell='E45'; Value='  -2.15%  '},
    @{Cell='D46'; Value='1.293'},
    @{Cell='E46'; Value='  -6.23%  '},
    @{Cell='D47'; Value='3.718'},
    @{Cell='E47'; Value='  -1.46%  '},
    @{Cell='D48'; Value='124.87'},
    @{Cell='E48'; Value='  -1.02%  '},
    @{Cell='D49'; Value='1.196'},
    @{Cell='E49'; Value='  -1.31%  '},
    @{Cell='D50'; Value='1.908'},
    @{Cell='E50'; Value='  -1.16%  '},
    @{Cell='D51'; Value='0.06731'}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
